$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$rng.Find.Replacement.Font.Superscript = $true
$found = $rng.Find.Execute([char]2, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
Write-Host "Found:" $found
